$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B2 (CasesTab query): append ORDER BY / LIMIT clause ---
$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100 "
$ws.Rows.Item(2).RowHeight = 345.6

# --- Update B3 (SamplesTab query): append ORDER BY / LIMIT clause ---
$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $b3 + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Rows.Item(3).RowHeight = 360

# --- Update B4 (FilesTab query): replace trailing "order by" with "order By ... LIMIT 100" ---
$b4 = $ws.Range("B4").Value2
$b4 = $b4 -replace "order by f\.file_name$", "order By f.file_name ASC LIMIT 100"
$ws.Range("B4").Value2 = $b4

# --- Update active selection / scroll position on the sheet ---
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 3
